{"js": "// Replace each three-digit-by-one-digit multiplication equation in the\n// document's table with its new value, per the commit's regenerated\n// answer key. Every original equation string is unique within the\n// document, so a simple exact-text search/replace per pair is safe.\nconst replacements = [\n  [\"874\u00d73=2622\", \"779\u00d76=4674\"],\n  [\"722\u00d79=6498\", \"293\u00d75=1465\"],\n  [\"256\u00d74=1024\", \"790\u00d79=7110\"],\n  [\"590\u00d75=2950\", \"423\u00d72=846\"],\n  [\"587\u00d78=4696\", \"854\u00d72=1708\"],\n  [\"245\u00d77=1715\", \"833\u00d75=4165\"],\n  [\"542\u00d74=2168\", \"548\u00d73=1644\"],\n  [\"778\u00d78=6224\", \"666\u00d76=3996\"],\n  [\"435\u00d77=3045\", \"591\u00d78=4728\"],\n  [\"878\u00d77=6146\", \"975\u00d72=1950\"],\n  [\"625\u00d73=1875\", \"614\u00d74=2456\"],\n  [\"581\u00d74=2324\", \"867\u00d79=7803\"],\n  [\"452\u00d77=3164\", \"285\u00d78=2280\"],\n  [\"552\u00d75=2760\", \"625\u00d77=4375\"],\n  [\"247\u00d77=1729\", \"117\u00d73=351\"],\n  [\"851\u00d74=3404\", \"158\u00d74=632\"],\n  [\"154\u00d78=1232\", \"779\u00d75=3895\"],\n  [\"159\u00d75=795\", \"306\u00d73=918\"],\n  [\"187\u00d76=1122\", \"479\u00d76=2874\"],\n  [\"421\u00d73=1263\", \"581\u00d72=1162\"],\n  [\"252\u00d74=1008\", \"675\u00d79=6075\"],\n  [\"124\u00d72=248\", \"209\u00d74=836\"],\n  [\"955\u00d74=3820\", \"242\u00d73=726\"],\n  [\"213\u00d73=639\", \"812\u00d78=6496\"],\n  [\"277\u00d75=1385\", \"992\u00d76=5952\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each three-digit-by-one-digit multiplication equation in the\n# document's table with its new value, per the commit's regenerated\n# answer key. Every original equation string is unique within the\n# document, so a simple Find/Replace (ReplaceAll) per pair is safe.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"874\u00d73=2622\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"779\u00d76=4674\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"722\u00d79=6498\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"293\u00d75=1465\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"256\u00d74=1024\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"790\u00d79=7110\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"590\u00d75=2950\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"423\u00d72=846\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"587\u00d78=4696\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"854\u00d72=1708\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"245\u00d77=1715\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"833\u00d75=4165\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"542\u00d74=2168\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"548\u00d73=1644\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"778\u00d78=6224\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"666\u00d76=3996\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"435\u00d77=3045\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"591\u00d78=4728\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"878\u00d77=6146\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"975\u00d72=1950\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"625\u00d73=1875\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"614\u00d74=2456\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"581\u00d74=2324\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"867\u00d79=7803\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"452\u00d77=3164\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"285\u00d78=2280\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"552\u00d75=2760\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"625\u00d77=4375\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"247\u00d77=1729\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"117\u00d73=351\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"851\u00d74=3404\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"158\u00d74=632\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"154\u00d78=1232\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"779\u00d75=3895\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"159\u00d75=795\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"306\u00d73=918\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"187\u00d76=1122\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"479\u00d76=2874\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"421\u00d73=1263\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"581\u00d72=1162\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"252\u00d74=1008\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"675\u00d79=6075\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"124\u00d72=248\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"209\u00d74=836\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"955\u00d74=3820\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"242\u00d73=726\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"213\u00d73=639\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"812\u00d78=6496\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"277\u00d75=1385\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"992\u00d76=5952\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n"}
